$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The test case name in A2 was renumbered/relabelled and the pickup /
# drop-off dates in C2 / D2 were bumped forward two months (xpath fix +
# data refresh per the commit message).
$ws.Range("A2").Value = "TC_02"
$ws.Range("C2").Value = "12/10/2017"
$ws.Range("D2").Value = "12/24/2017"

# The rest of the used range (everything except the already explicitly
# styled header cell A1) had its font explicitly set to Arial 10 -
# matching the workbook's default font - mirroring the formatting pass
# that produced the duplicated-but-equivalent font entries in the saved
# file.
$ws.Range("B1:H1").Font.Name = "Arial"
$ws.Range("B1:H1").Font.Size = 10
$ws.Range("A2:H2").Font.Name = "Arial"
$ws.Range("A2:H2").Font.Size = 10

# Active cell moved to D2 (the updated pickup date) before saving.
$ws.Range("D2").Select() | Out-Null
